$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns G and H (TOC Coverage (%) and Section Coverage (%))
$ws.Range("G1:H2").EntireColumn.Delete()

# Update changed values in row 2
$ws.Range("B2").Value = 922
$ws.Range("C2").Value = 3635
$ws.Range("E2").Value = 1014
